$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136; this shifts the existing rows 136-145
# down to 137-146, preserving their values and formatting.
$ws.Rows("136:136").Insert()

# Copy the date cell's number formatting from the row above (row 135) so the
# new date cell in D136 renders like the rest of the column.
$ws.Range("D135").Copy()
$ws.Range("D136").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's data (same static columns as surrounding rows, with
# the new record's values).
$ws.Range("A136").Value = 6
$ws.Range("B136").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C136").Value = "Metropolitana"
$ws.Range("D136").Value = 44578
$ws.Range("E136").Value = 13
$ws.Range("F136").Value = 100112029
$ws.Range("G136").Value = "Orégano"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 32
$ws.Range("K136").Value = 8500
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 8719
$ws.Range("N136").Value = "$/docena de atados"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 2906
$ws.Range("Q136").Value = 3
$ws.Range("R136").Value = "Hortaliza"
